$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.975.81'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '2.269.77'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.17'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.629'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.69'
$ws.Range("E7").Value = '  +3.73%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.448'
$ws.Range("E9").Value = '  +6.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0995'
$ws.Range("E10").Value = '  +5.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.48'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.31'
$ws.Range("E12").Value = '  +14.76%  '
$ws.Range("E13").Value = '  +2.06%  '
$ws.Range("D14").Value = '2.611.18'
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.79'
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.12'
$ws.Range("E16").Value = '  +5.78%  '
$ws.Range("D18").Value = '2.280.08'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").Value = '43.926.25'
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000101'
$ws.Range("E20").Value = '  +7.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.82'
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.64'
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("E25").Value = '  -3.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  -1.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.33'
$ws.Range("E27").Value = '  +25.65%  '
$ws.Range("E28").Value = '  +2.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.47'
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("E30").Value = '  -1.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.92'
$ws.Range("E31").Value = '  +1.76%  '
$ws.Range("E32").Value = '  -5.30%  '
$ws.Range("E33").Value = '  +2.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0706'
$ws.Range("E34").Value = '  +6.77%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.86'
$ws.Range("E36").Value = '  -3.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.82'
$ws.Range("E37").Value = '  +4.80%  '
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("E39").Value = '  -3.69%  '
$ws.Range("E40").Value = '  +3.14%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000225'
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0989'
$ws.Range("E43").Value = '  +2.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.56'
$ws.Range("E44").Value = '  +4.51%  '
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.50'
$ws.Range("E45").Value = '  +9.80%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.27'
$ws.Range("E46").Value = '  -5.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.26'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.39'
$ws.Range("E49").Value = '  -1.83%  '
$ws.Range("D50").Value = '1.446.69'
$ws.Range("E50").Value = '  -1.69%  '
$ws.Range("E51").Value = '  +1.43%  '
